$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the
#    first (title) paragraph.
# ---------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$boldStart = $metaRange.Start
$metaRange.InsertAfter("Meta description")
$boldEnd = $boldStart + ("Meta description").Length

$metaPara2 = $d.Paragraphs.Item(2)
$tailRange = $metaPara2.Range
$tailRange.InsertAfter(": Experience the Rising Re-Spins feature in Apollo Rising slot game free. RTP analysis, comparisons, and potential jackpot winnings explained in our review.")

$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Font.Bold = 1

# ---------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of
#    the document ("Play Apollo Rising Slot Game Free - Review &
#    RTP Analysis"). After step 1 the document grew by one
#    paragraph, so the paragraph that used to be #47 is now #48.
# ---------------------------------------------------------------
$dupTitlePara = $d.Paragraphs.Item(48)
$dupTitlePara.Range.Delete()

# ---------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the
#    DALL-E image prompt, keeping its italic formatting. After
#    the delete above it is paragraph #48 again (the last one).
# ---------------------------------------------------------------
$blurbPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$blurbRange = $blurbPara.Range.Duplicate()
$blurbRange.Find.Execute("Experience the Rising Re-Spins feature in Apollo Rising slot game free. RTP analysis, comparisons, and potential jackpot winnings explained in our review.") | Out-Null
$blurbRange.Text = 'DALLE, please create a feature image for the game "Apollo Rising". The image should be in cartoon style and feature a happy Maya warrior with glasses. This should fit the space theme of the game, with the background including rich graphics of neon blue shades, stars, and rockets to create an atmosphere of a space mission. The image should capture the exciting and innovative gameplay of the game while also incorporating the Maya warrior with glasses to add a unique touch. Please make it eye-catching and attention-grabbing to draw in potential players. Thank you!'
